# Apply cell-value updates captured in the "Updated symbol list" data refresh.
# Columns D (Price) and E (Volume(1h)) hold numeric/percentage-looking text; a leading
# apostrophe forces Excel to keep them as literal text (matching the source t="inlineStr"
# cells) instead of silently coercing them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '''327.86'
    'E2' = '''-0.54%'
    'D3' = '''43.83'
    'E3' = '''0.71%'
    'D4' = '''5.538'
    'E4' = '''-1.17%'
    'D5' = '''0.08016'
    'E5' = '''-2.27%'
    'D6' = '''1.894'
    'E6' = '''0.54%'
    'D7' = '''4.258'
    'E7' = '''-2.45%'
    'B8' = 'BTSEToken'
    'C8' = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
    'D8' = '''2.560'
    'E8' = '''-8.72%'
    'B9' = 'MXToken'
    'C9' = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    'D9' = '''0.9430'
    'E9' = '''-0.05%'
    'D10' = '''0.1188'
    'E10' = '''-0.24%'
    'E11' = '''-3.78%'
    'B12' = 'MandalaExchangeToken'
    'C12' = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
    'D12' = '''0.09643'
    'E12' = '''-1.09%'
    'B13' = 'BitrueCoin'
    'C13' = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
    'D13' = '''0.04435'
    'E13' = '''2.62%'
    'B14' = 'BitMartToken'
    'C14' = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
    'D14' = '''0.1066'
    'E14' = '''-0.38%'
    'B15' = 'BitForexToken'
    'C15' = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
    'D15' = '''0.001273'
    'E15' = '''-0.67%'
    'B16' = 'TigerCash'
    'C16' = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
    'D16' = '''0.005938'
    'E16' = '''-0.03%'
    'B17' = 'LEO'
    'C17' = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
    'D17' = '''3.405'
    'E17' = '''-3.56%'
    'B18' = 'BitpandaEcosystemToken'
    'C18' = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
    'D18' = '''0.3428'
    'E18' = '''-3.06%'
    'B19' = 'MCDex'
    'C19' = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
    'D19' = '''10.05'
    'E19' = '''14.98%'
    'D20' = '''0.1414'
    'E20' = '''3.22%'
    'D22' = '''0.04194'
    'E22' = '''-4.59%'
    'D23' = '''0.001247'
    'E23' = '''0.46%'
    'D24' = '''0.004283'
    'E24' = '''-0.42%'
    'D25' = '''0.0001262'
    'E25' = '''2.05%'
    'D26' = '''0.0003992'
    'E26' = '''-0.46%'
    'D38' = '''0.02642'
    'E38' = '''-3.93%'
    'D39' = '''0.05485'
    'E39' = '''-3.45%'
    'D40' = '''0.007579'
    'E40' = '''-3.98%'
    'E41' = '''-2.03%'
    'D42' = '''0.008138'
    'E42' = '''-16.62%'
    'D43' = '''0.002003'
    'E43' = '''-4.98%'
    'D44' = '''0.008799'
    'E44' = '''-12.26%'
    'D45' = '''0.00007102'
    'E45' = '''-3.04%'
    'D46' = '''0.00000000751'
    'E46' = '''-0.45%'
    'B47' = 'BOLO'
    'C47' = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
    'D47' = '''0.003278'
    'E47' = '''-5.14%'
    'B48' = 'CoinbaseStockToken'
    'C48' = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
    'D48' = '''0.002272'
    'E48' = '''-0.46%'
    'D49' = '''0.00002102'
    'E49' = '''-0.45%'
    'D50' = '''0.0002002'
    'E50' = '''-0.45%'
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
